# Update crypto price/volume figures for the Mon Apr 24 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "27.985.78"
$ws.Range("D3").Value = "1.879.40"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").Value = "'336.43"
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").Value = "'0.3948"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'47.09"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "'0.08016"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "1.885.40"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "'6.056"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "'7.203"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'1.014"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "'88.71"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "'0.06744"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "'1.010"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "27.977.98"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "'5.510"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'11.02"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "'2.343"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "2.108.56"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'158.61"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'19.89"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "'2.106"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "'5.489"
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").Value = "'121.56"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "'0.9789"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'0.09554"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "'5.343"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'1.355"
$ws.Range("E36").Value = "  -6.54%  "
$ws.Range("D37").Value = "'0.06088"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "'0.02246"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'8.209"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'1.010"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "'0.5987"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "'0.1899"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'10.35"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "'1.257"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'0.5682"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "'12.21"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "'3.342"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").Value = "'0.06787"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "'112.70"
$ws.Range("E51").Value = "  -1.68%  "
